$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold font, border, centered) from H1 to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set header values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Set data values
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
